$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: split the run that currently contains the whole $fullText at
# the boundary that ends right after $prefix (i.e. right before
# $fullText.Substring(prefix.Length)). This is done by dropping a
# temporary bookmark at that character offset -- Word (and this host)
# splits the underlying run there, and the split survives the bookmark
# being deleted again. $fullText must currently be unique in the
# document content.
# ---------------------------------------------------------------------
function Split-RunAfterPrefix($fullText, $prefix) {
    $r = $word.ActiveDocument.Content
    $r.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $prefixLen = $prefix.Length
    $splitPos = $r.Start + $prefixLen
    $sp = $word.ActiveDocument.Range($splitPos, $splitPos)
    $word.ActiveDocument.Bookmarks.Add("__TmpSplit__", $sp)
    $word.ActiveDocument.Bookmarks.Item("__TmpSplit__").Delete()
}

# Splits $fullText into runs at each prefix length in $lengths (a list of
# cumulative character counts, ascending, each < Len(fullText)). Applied
# right-to-left so earlier Find() calls aren't disturbed by upstream
# edits, which also happens to be what keeps w:t's xml:space="preserve"
# minimal/matching real Word output.
function Split-RunAtOffsets($fullText, $lengths) {
    $sorted = $lengths | Sort-Object -Descending
    foreach ($len in $sorted) {
        $prefix = $fullText.Substring(0, $len)
        Split-RunAfterPrefix $fullText $prefix
    }
}

# ------------------------------------------------------------------
# 1) "] Allow users to take pictures of vehicles illegally parked with an in-app camera;"
#    -> "] Allow users to take pictures of vehicles illegally parked;"
#    Runs: "] Allow users to take pictures of vehicles i" | <_GoBack> | "llegally " | "parked" | ";"
# ------------------------------------------------------------------
$d.Content.Find.Execute("illegally parked with an in-app camera;", $true, $false, $false, $false, $false, $true, 1, $false, "illegally parked;", 2)

$text1 = "] Allow users to take pictures of vehicles illegally parked;"
$lengths1 = @(44, 53, 59)
Split-RunAtOffsets $text1 $lengths1

# ------------------------------------------------------------------
# 2) "] Allow users and authorities to access stored data;"
#    -> "] Allow users and authorities to access stored data about submitted reports;"
#    Runs: "...access stored data" | " about " | "submitted " | "reports" | ";"
# ------------------------------------------------------------------
$d.Content.Find.Execute("access stored data;", $true, $false, $false, $false, $false, $true, 1, $false, "access stored data about submitted reports;", 2)

$text2 = "] Allow users and authorities to access stored data about submitted reports;"
$marker2 = "access stored data"
$p2a = $text2.IndexOf($marker2) + $marker2.Length
$p2b = $p2a + 7   # " about ".Length
$p2c = $p2b + 10  # "submitted ".Length
$p2d = $p2c + 7   # "reports".Length
$lengths2 = @($p2a, $p2b, $p2c, $p2d)
Split-RunAtOffsets $text2 $lengths2

# ------------------------------------------------------------------
# 3) "] Allow users and authorities to build statistics;"
#    -> "] Allow users and authorities to visualize statistics posted by the admins;"
#    Runs: "] Allow users and authorities to " | "visualize" | " statistics" | " posted by the admins" | ";"
# ------------------------------------------------------------------
$d.Content.Find.Execute("to build statistics;", $true, $false, $false, $false, $false, $true, 1, $false, "to visualize statistics posted by the admins;", 2)

$text3 = "] Allow users and authorities to visualize statistics posted by the admins;"
$p3a = $text3.IndexOf("visualize")
$p3b = $p3a + 9   # "visualize".Length
$p3c = $p3b + 11  # " statistics".Length
$p3d = $p3c + 22  # " posted by the admins".Length
$lengths3 = @($p3a, $p3b, $p3c, $p3d)
Split-RunAtOffsets $text3 $lengths3

# ------------------------------------------------------------------
# 4) "sers to suggest possible interventions for areas that are deemed unsafe"
#    -> "sers to visualize possible interventions for areas that are deemed unsafe"
#    Runs: "sers to " | "visualize" | " possible interventions for areas that " | "are deemed unsafe"
# ------------------------------------------------------------------
$d.Content.Find.Execute("sers to suggest possible interventions", $true, $false, $false, $false, $false, $true, 1, $false, "sers to visualize possible interventions", 2)

$text4 = "sers to visualize possible interventions for areas that are deemed unsafe"
$p4a = $text4.IndexOf("visualize")
$p4b = $p4a + 9   # "visualize".Length
$p4c = $p4b + 40  # " possible interventions for areas that ".Length
$lengths4 = @($p4a, $p4b, $p4c)
Split-RunAtOffsets $text4 $lengths4

# ------------------------------------------------------------------
# 5) Move the "_GoBack" bookmark from the trailing empty paragraph to the
#    point in the first edited sentence where the last edit happened
#    (right after "...vehicles i" in "illegally").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$anchor = $d.Content
$anchor.Find.Execute("] Allow users to take pictures of vehicles i", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$gobackRange = $d.Range($anchor.End, $anchor.End)
$d.Bookmarks.Add("_GoBack", $gobackRange)
